# docs/Data Feed Specification.xlsx
# "chinh sua datafeed, them thiet ke warehouse"
#
# The "Example" column for the first two fields (product_id / product_name)
# in both mini-tables (PNJ: F3:F4, Dang Quang Watch: M3:M4) had their sample
# values swapped - product_id's example should be the SKU-looking string and
# product_name's example should be the descriptive product name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap the "Example" values for product_id / product_name rows ---
$f3 = $ws.Range("F3").Value()
$f4 = $ws.Range("F4").Value()
$m3 = $ws.Range("M3").Value()
$m4 = $ws.Range("M4").Value()

$ws.Range("F3").Value = $f4
$ws.Range("F4").Value = $f3
$ws.Range("M3").Value = $m4
$ws.Range("M4").Value = $m3

# --- restyle: row 3 (product_id example) gets a plain black font, ---
# --- row 4 (product_name example) drops back to the column's default (no border) ---
$ws.Range("F3").Font.Color = 0
$ws.Range("M3").Font.Color = 0
$ws.Range("F4").Borders.LineStyle = -4142
$ws.Range("M4").Borders.LineStyle = -4142

# --- move the saved selection cursor ---
$null = $ws.Range("F20").Select()
